$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2.4
$ws.Range("G2").Value = 3.4
$ws.Range("H2").Value = 2.52
$ws.Range("I2").Value = 3.6
$ws.Range("J2").Value = 2.84
$ws.Range("K2").Value = 3.9
$ws.Range("M2").Value = 1.09
$ws.Range("N2").Value = 2.52
$ws.Range("O2").Value = 1.41
$ws.Range("P2").Value = 1.57
$ws.Range("Q2").Value = 2.1
$ws.Range("R2").Value = 1.21
$ws.Range("S2").Value = 3.45
$ws.Range("T2").Value = 1.86
$ws.Range("U2").Value = 1.81
$ws.Range("V2").Value = 1.39
$ws.Range("W2").Value = 1.42
$ws.Range("F3").Value = 1.95
$ws.Range("G3").Value = 2.36
$ws.Range("H3").Value = 3.7
$ws.Range("I3").Value = 5.4
$ws.Range("J3").Value = 2.9
$ws.Range("K3").Value = 3.7
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 2.42
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 1.53
$ws.Range("Q3").Value = 2.2
$ws.Range("R3").Value = 1.2
$ws.Range("S3").Value = 3.45
$ws.Range("T3").Value = 1.98
$ws.Range("U3").Value = 1.72
$ws.Range("V3").Value = 1.23
$ws.Range("W3").Value = 1.73
$ws.Range("X3").Value = 11.5
$ws.Range("Y3").Value = 13.5
$ws.Range("Z3").Value = 36
$ws.Range("AB3").Value = 8
$ws.Range("AC3").Value = 8.199999999999999
$ws.Range("AD3").Value = 21
$ws.Range("AF3").Value = 13.5
$ws.Range("AG3").Value = 12.5
$ws.Range("AH3").Value = 26
$ws.Range("AJ3").Value = 32
$ws.Range("AK3").Value = 32
$ws.Range("AN3").Value = 32
$ws.Range("F4").Value = 5.9
$ws.Range("G4").Value = 6
$ws.Range("H4").Value = 1.74
$ws.Range("I4").Value = 1.89
$ws.Range("J4").Value = 3.2
$ws.Range("L4").Value = 1.6
$ws.Range("N4").Value = 2.38
$ws.Range("O4").Value = 1.6
$ws.Range("P4").Value = 1.45
$ws.Range("Q4").Value = 2.76
$ws.Range("R4").Value = 1.16
$ws.Range("S4").Value = 6
$ws.Range("T4").Value = 2.44
$ws.Range("U4").Value = 1.57
$ws.Range("V4").Value = 2.12
$ws.Range("W4").Value = 1.2
$ws.Range("X4").Value = 8.199999999999999
$ws.Range("Y4").Value = 5.7
$ws.Range("Z4").Value = 9
$ws.Range("AA4").Value = 22
$ws.Range("AB4").Value = 15
$ws.Range("AC4").Value = 8.6
$ws.Range("AD4").Value = 12
$ws.Range("AE4").Value = 29
$ws.Range("AF4").Value = 55
$ws.Range("AG4").Value = 30
$ws.Range("AH4").Value = 44
$ws.Range("AI4").Value = 90
$ws.Range("AJ4").Value = 1000
$ws.Range("AK4").Value = 180
$ws.Range("AL4").Value = 1000
$ws.Range("AM4").Value = 440
$ws.Range("AN4").Value = 1000
$ws.Range("AO4").Value = 24
$ws.Range("G5").Value = 4.5
$ws.Range("H5").Value = 2.12
$ws.Range("I5").Value = 2.46
$ws.Range("N5").Value = 2.56
$ws.Range("R5").Value = 1.19
$ws.Range("V5").Value = 1.68
$ws.Range("L6").Value = 1.53
$ws.Range("G7").Value = 3.5
$ws.Range("H7").Value = 2.78
$ws.Range("I7").Value = 2.98
$ws.Range("J7").Value = 2.72
$ws.Range("O7").Value = 1.83
$ws.Range("P7").Value = 1.35
$ws.Range("V7").Value = 1.51
$ws.Range("Y7").Value = 7.2
$ws.Range("Z7").Value = 15.5
$ws.Range("AB7").Value = 7.8
$ws.Range("AC7").Value = 7.4
$ws.Range("AG7").Value = 18.5
$ws.Range("AI7").Value = 120
$ws.Range("AN7").Value = 1000
$ws.Range("F8").Value = 2.2
$ws.Range("H8").Value = 3.8
$ws.Range("J8").Value = 2.96
$ws.Range("L8").Value = 1.59
$ws.Range("N8").Value = 2.54
$ws.Range("O8").Value = 1.56
$ws.Range("Q8").Value = 2.6
$ws.Range("S8").Value = 5.6
$ws.Range("U8").Value = 1.71
$ws.Range("X8").Value = 980
$ws.Range("F9").Value = 2
$ws.Range("G9").Value = 2.12
$ws.Range("H9").Value = 4.5
$ws.Range("K9").Value = 3.4
$ws.Range("F10").Value = 2.4
$ws.Range("H10").Value = 3.4
$ws.Range("I10").Value = 3.6
$ws.Range("J10").Value = 3.15
$ws.Range("N10").Value = 2.5
$ws.Range("P10").Value = 1.5
$ws.Range("Q10").Value = 2.76
$ws.Range("R10").Value = 1.16
$ws.Range("T10").Value = 2.28
$ws.Range("U10").Value = 1.7
$ws.Range("V10").Value = 1.39
$ws.Range("Y10").Value = 8.800000000000001
$ws.Range("Z10").Value = 22
$ws.Range("AA10").Value = 80
$ws.Range("AC10").Value = 7.4
$ws.Range("AD10").Value = 21
$ws.Range("AG10").Value = 15
$ws.Range("AI10").Value = 95
$ws.Range("AJ10").Value = 38
$ws.Range("AK10").Value = 75
